# تعديل تلقائي في شيت Card1 by admin at 2025-11-05 19:13:24
# Update row 6 of the "Card1" sheet:
#  - Fill the previously-empty D6:K6 and M6 cells with the literal text "nan"
#  - Change the date in L6 from "29\9\2024" to "22\12\2024"

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card1")

$ws.Range("D6").Value = "nan"
$ws.Range("E6").Value = "nan"
$ws.Range("F6").Value = "nan"
$ws.Range("G6").Value = "nan"
$ws.Range("H6").Value = "nan"
$ws.Range("I6").Value = "nan"
$ws.Range("J6").Value = "nan"
$ws.Range("K6").Value = "nan"
$ws.Range("L6").Value = "22\12\2024"
$ws.Range("M6").Value = "nan"
